# Nieuwe data toegevoegd via Streamlit op 2024-12-04 11:20:37
# Append a new row (row 94) with the new KDV entry to Sheet1.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newRow = 94

# Columns A-C: plain text values
$ws.Cells.Item($newRow, 1).Value = "Kindergarden"
$ws.Cells.Item($newRow, 2).Value = "Kindergarden Delft Buitenhofdreef"
$ws.Cells.Item($newRow, 3).Value = "KDV"

# Column D: date stored as literal text (not an Excel date serial number),
# matching the plain "yyyy-mm-dd" format used by the most recently added rows.
# Force a text number format while assigning so Excel does not coerce the
# string into a date value, then restore the default cell style so the
# cell keeps looking like the rest of the (unstyled) data rows.
$ws.Cells.Item($newRow, 4).NumberFormat = "@"
$ws.Cells.Item($newRow, 4).Value = "2023-12-07"
$ws.Cells.Item($newRow, 4).Style = "Normal"

# Columns E-J: numeric zero values
$ws.Cells.Item($newRow, 5).Value = 0
$ws.Cells.Item($newRow, 6).Value = 0
$ws.Cells.Item($newRow, 7).Value = 0
$ws.Cells.Item($newRow, 8).Value = 0
$ws.Cells.Item($newRow, 9).Value = 0
$ws.Cells.Item($newRow, 10).Value = 0
